$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header column O1 = "assignsurveycoordinator" (appends to the
# shared strings table), copying the formatting (style) already used by
# the preceding header cell N1.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("O1").Value = "assignsurveycoordinator"

# Row 1 grows to fit the now 3-line wrapped header text.
$ws.Rows.Item(1).RowHeight = 43.5

# Update the selected cell to K3, matching the saved selection in the file.
$ws.Range("K3").Select()
